{"js": "const replacements = [\n  [\"85\u00d742=3570\", \"16\u00d739=624\"],\n  [\"99\u00d711=1089\", \"57\u00d718=1026\"],\n  [\"92\u00d753=4876\", \"32\u00d713=416\"],\n  [\"89\u00d723=2047\", \"78\u00d721=1638\"],\n  [\"42\u00d718=756\", \"60\u00d749=2940\"],\n  [\"92\u00d725=2300\", \"38\u00d715=570\"],\n  [\"44\u00d721=924\", \"93\u00d776=7068\"],\n  [\"14\u00d793=1302\", \"79\u00d764=5056\"],\n  [\"91\u00d776=6916\", \"87\u00d715=1305\"],\n  [\"75\u00d740=3000\", \"56\u00d716=896\"],\n  [\"86\u00d797=8342\", \"60\u00d736=2160\"],\n  [\"26\u00d752=1352\", \"50\u00d798=4900\"],\n  [\"31\u00d766=2046\", \"48\u00d717=816\"],\n  [\"25\u00d736=900\", \"62\u00d732=1984\"],\n  [\"60\u00d788=5280\", \"46\u00d793=4278\"],\n  [\"84\u00d712=1008\", \"26\u00d749=1274\"],\n  [\"91\u00d745=4095\", \"78\u00d767=5226\"],\n  [\"41\u00d736=1476\", \"97\u00d746=4462\"],\n  [\"67\u00d784=5628\", \"48\u00d789=4272\"],\n  [\"82\u00d762=5084\", \"98\u00d795=9310\"],\n  [\"54\u00d770=3780\", \"20\u00d739=780\"],\n  [\"11\u00d799=1089\", \"58\u00d797=5626\"],\n  [\"83\u00d794=7802\", \"48\u00d789=4272\"],\n  [\"59\u00d717=1003\", \"73\u00d769=5037\"],\n  [\"31\u00d786=2666\", \"42\u00d779=3318\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst map = new Map(replacements);\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (map.has(text)) {\n    paragraph.insertText(map.get(text), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"85\u00d742=3570\"; Replace = \"16\u00d739=624\" },\n    @{ Find = \"99\u00d711=1089\"; Replace = \"57\u00d718=1026\" },\n    @{ Find = \"92\u00d753=4876\"; Replace = \"32\u00d713=416\" },\n    @{ Find = \"89\u00d723=2047\"; Replace = \"78\u00d721=1638\" },\n    @{ Find = \"42\u00d718=756\"; Replace = \"60\u00d749=2940\" },\n    @{ Find = \"92\u00d725=2300\"; Replace = \"38\u00d715=570\" },\n    @{ Find = \"44\u00d721=924\"; Replace = \"93\u00d776=7068\" },\n    @{ Find = \"14\u00d793=1302\"; Replace = \"79\u00d764=5056\" },\n    @{ Find = \"91\u00d776=6916\"; Replace = \"87\u00d715=1305\" },\n    @{ Find = \"75\u00d740=3000\"; Replace = \"56\u00d716=896\" },\n    @{ Find = \"86\u00d797=8342\"; Replace = \"60\u00d736=2160\" },\n    @{ Find = \"26\u00d752=1352\"; Replace = \"50\u00d798=4900\" },\n    @{ Find = \"31\u00d766=2046\"; Replace = \"48\u00d717=816\" },\n    @{ Find = \"25\u00d736=900\"; Replace = \"62\u00d732=1984\" },\n    @{ Find = \"60\u00d788=5280\"; Replace = \"46\u00d793=4278\" },\n    @{ Find = \"84\u00d712=1008\"; Replace = \"26\u00d749=1274\" },\n    @{ Find = \"91\u00d745=4095\"; Replace = \"78\u00d767=5226\" },\n    @{ Find = \"41\u00d736=1476\"; Replace = \"97\u00d746=4462\" },\n    @{ Find = \"67\u00d784=5628\"; Replace = \"48\u00d789=4272\" },\n    @{ Find = \"82\u00d762=5084\"; Replace = \"98\u00d795=9310\" },\n    @{ Find = \"54\u00d770=3780\"; Replace = \"20\u00d739=780\" },\n    @{ Find = \"11\u00d799=1089\"; Replace = \"58\u00d797=5626\" },\n    @{ Find = \"83\u00d794=7802\"; Replace = \"48\u00d789=4272\" },\n    @{ Find = \"59\u00d717=1003\"; Replace = \"73\u00d769=5037\" },\n    @{ Find = \"31\u00d786=2666\"; Replace = \"42\u00d779=3318\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
